$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A6: date value changes from 45295 to 45294 (keep underlying numeric/date typing)
$ws.Range("A6").Value = 45294

# E8, E10, E12: truck_used counts change from 1 to 2
$ws.Range("E8").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("E12").Value = 2
